$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp update (A1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 7 de Septiembre de 2020 a las 02:09"

# Row 4 - Estados Unidos
$ws.Cells.Item(4,2).Value = 6459346
$ws.Cells.Item(4,3).Value = 30206
$ws.Cells.Item(4,4).Value = 3725476
$ws.Cells.Item(4,5).Value = 2540634
$ws.Cells.Item(4,7).Value = 418
$ws.Cells.Item(4,8).Value = 193236

# Row 6 - Brasil
$ws.Cells.Item(6,2).Value = 4137606
$ws.Cells.Item(6,3).Value = 14606
$ws.Cells.Item(6,5).Value = 693693
$ws.Cells.Item(6,7).Value = 456
$ws.Cells.Item(6,8).Value = 126686

# Row 8 - Peru
$ws.Cells.Item(8,2).Value = 689977
$ws.Cells.Item(8,3).Value = 6275
$ws.Cells.Item(8,4).Value = 515039
$ws.Cells.Item(8,5).Value = 145100
$ws.Cells.Item(8,7).Value = 151
$ws.Cells.Item(8,8).Value = 29838

# Row 24 - Alemania
$ws.Cells.Item(24,2).Value = 251724
$ws.Cells.Item(24,3).Value = 668
$ws.Cells.Item(24,5).Value = 16115

# Row 28 - Canada
$ws.Cells.Item(28,2).Value = 131895
$ws.Cells.Item(28,3).Value = 400
$ws.Cells.Item(28,4).Value = 116357
$ws.Cells.Item(28,5).Value = 6393

# Row 36 - Panama
$ws.Cells.Item(36,2).Value = 97043
$ws.Cells.Item(36,3).Value = 738
$ws.Cells.Item(36,4).Value = 69661
$ws.Cells.Item(36,5).Value = 25296
$ws.Cells.Item(36,7).Value = 11
$ws.Cells.Item(36,8).Value = 2086

# Row 72 - Chequia
$ws.Cells.Item(72,2).Value = 28156
$ws.Cells.Item(72,3).Value = 404
$ws.Cells.Item(72,4).Value = 19053
$ws.Cells.Item(72,5).Value = 8667
$ws.Cells.Item(72,7).Value = 5
$ws.Cells.Item(72,8).Value = 436

# Row 106 - Luxemburgo
$ws.Cells.Item(106,2).Value = 6950
$ws.Cells.Item(106,3).Value = 54
$ws.Cells.Item(106,5).Value = 700

# Row 109 - Montenegro
$ws.Cells.Item(109,2).Value = 5553
$ws.Cells.Item(109,3).Value = 131
$ws.Cells.Item(109,4).Value = 4231
$ws.Cells.Item(109,5).Value = 1214
$ws.Cells.Item(109,7).Value = 1
$ws.Cells.Item(109,8).Value = 108

# Rows 121/122 - Cabo Verde and Surinam swap order (Surinam now sorts before
# Cabo Verde) and Surinam's stats are refreshed; Cabo Verde keeps its old
# values, just shifted down one row.
$ws.Cells.Item(121,1).Value = "Surinam"
$ws.Cells.Item(121,2).Value = 4346
$ws.Cells.Item(121,3).Value = 26
$ws.Cells.Item(121,4).Value = 3494
$ws.Cells.Item(121,5).Value = 767
$ws.Cells.Item(121,6).Value = 0
$ws.Cells.Item(121,7).Value = 9
$ws.Cells.Item(121,8).Value = 85

$ws.Cells.Item(122,1).Value = "Cabo Verde"
$ws.Cells.Item(122,2).Value = 4330
$ws.Cells.Item(122,3).Value = 55
$ws.Cells.Item(122,4).Value = 3628
$ws.Cells.Item(122,5).Value = 660
$ws.Cells.Item(122,6).Value = 0
$ws.Cells.Item(122,7).Value = 0
$ws.Cells.Item(122,8).Value = 42
